# appendix-z.docx — "Updates to Chef DK 0.17.17 and InSpec"
#
# The only substantive textual change in this revision is the version
# bump on the lab-image selection line:
#
#   "TDD Cookbook Development [EN DASH] CentOS 6.7 [EN DASH] 1.0.0"
#     -> "TDD Cookbook Development [EN DASH] CentOS 6.7 [EN DASH] 1.1.0"
#
# (The rest of the underlying OOXML diff only wraps "CentOS" and a
# couple of other already-correct spans in <w:proofErr> spell/grammar
# -check markers and splits the surrounding runs at those boundaries.
# Word's proofing engine inserts those automatically while it re-checks
# the document on edit/save; there is no accompanying text or
# formatting change, and <w:proofErr> is not a node the Word object
# model exposes for scripts to create, so there is nothing further for
# this script to do there.)

$d = $word.ActiveDocument

$enDash = [char]0x2013
$oldLine = "TDD Cookbook Development $enDash CentOS 6.7 $enDash 1.0.0"
$newLine = "TDD Cookbook Development $enDash CentOS 6.7 $enDash 1.1.0"

# Primary: replace the whole title line in one shot so the surrounding
# text/formatting is left completely untouched.
$found = $d.Content.Find.Execute($oldLine, $false, $true, $false, $false,
                                  $false, $true, 1, $false, $newLine, 2)

if (-not $found) {
    # Fallback (only needed if the line is already split across runs,
    # punctuation/spacing differs slightly, etc.): target just the
    # version number, which is unique in the document.
    $found = $d.Content.Find.Execute("1.0.0", $false, $true, $false, $false,
                                      $false, $true, 1, $false, "1.1.0", 2)
}

Write-Output "Version bump applied: $found"
